# TokenIteratorFieldRewriterSplit update: the M2Doc parser now rewrites
# the Word "fldChar begin / instrText / fldChar end" field triples that
# hold M2Doc template markers into plain literal text runs that read
# "{<expression>}" -- no more real Word fields, just braces typed as text.
#
# This document (extraSpaceInEndTemplate-template.docx) has four such
# fields, one per paragraph. For each field: take its field code (the
# M2Doc expression), strip the single leading/trailing space that
# surrounds it, wrap it in "{" / "}", delete the field and retype the
# paragraph as that literal text. The last field's code additionally
# contains a (hidden) "_GoBack" bookmark partway through; that bookmark
# is re-created at the same relative offset once the literal text is
# back in place.

$d = $word.ActiveDocument

function Get-ParagraphIndexForPosition($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $rng = $doc.Paragraphs.Item($i).Range
        if ($pos -ge $rng.Start -and $pos -lt $rng.End) {
            return $i
        }
    }
    return $doc.Paragraphs.Count
}

function ConvertTo-LiteralBraceText($code) {
    $body = $code
    if ($body.StartsWith(" ")) { $body = $body.Substring(1) }
    if ($body.EndsWith(" ")) { $body = $body.Substring(0, $body.Length - 1) }
    return "{" + $body + "}"
}

while ($d.Fields.Count -gt 0) {
    $f = $d.Fields.Item(1)
    $code = $f.Code.Text
    $paraIndex = Get-ParagraphIndexForPosition $d $f.Code.Start

    # The "endtemplate" field carries a hidden _GoBack bookmark right
    # after the "m: " that precedes it -- remember where, relative to
    # the field code, before the field is deleted.
    $hasGoBack = $code.Contains("endtemplate")
    $bookmarkOffsetInCode = 4  # "{m: " is 4 characters before the bookmark

    $literal = ConvertTo-LiteralBraceText $code
    $f.Delete()

    $p = $d.Paragraphs.Item($paraIndex).Range
    $p.Text = $literal
    $p.Font.LanguageID = "en-US"

    if ($hasGoBack) {
        $p2 = $d.Paragraphs.Item($paraIndex).Range
        $bookmarkPos = $p2.Start + $bookmarkOffsetInCode
        $bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
        $d.Bookmarks.Add("_GoBack", $bookmarkRange)
    }
}
